# Apply the updates made to TestCommittedProjects.xlsx:
#  - Committed Projects sheet, cell D2 (year) changes from 2021 to 2022
#  - The active selection on the sheet moves from B6 to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Committed Projects")

# Update the year value in D2
$ws.Range("D2").Value = 2022

# Make sure this is the active sheet and select D2 (matches saved selection in XML)
$ws.Activate()
$ws.Range("D2").Select()
